# Applies the scheduled-runner market data refresh to the Phantom_Profits workbook.
# For each affected row, updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N)
# to the freshly captured values; some rows gain or lose individual profit cells
# depending on whether NQ/HQ profit is present after the refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 2472.75
$ws.Range("I17").Value = 2995
$ws.Range("J17").Value = 2298.6667
$ws.Range("K17").Value = 8985
$ws.Range("L17").Value = 6896.000100000001
$ws.Range("M17").Value = -8817
$ws.Range("N17").Value = -7232.000100000001
# row 29
$ws.Range("H29").Value = 20499
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 25373.75
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 76121.25
$ws.Range("M29").Value = -2719
$ws.Range("N29").Value = -76683.25
# row 33
$ws.Range("H33").Value = 403
$ws.Range("I33").Value = 415.30768
$ws.Range("J33").Value = 349.66666
$ws.Range("K33").Value = 415.30768
$ws.Range("L33").Value = 349.66666
$ws.Range("M33").Value = -186.30768
$ws.Range("N33").Value = -807.66666
# row 38
$ws.Range("H38").Value = 291.1111
$ws.Range("I38").Value = 291.1111
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 873.3333
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -501.3333
$ws.Range("N38").ClearContents()
# row 58
$ws.Range("H58").Value = 2397
$ws.Range("I58").Value = 255
$ws.Range("J58").Value = 3825
$ws.Range("K58").Value = 765
$ws.Range("L58").Value = 11475
$ws.Range("M58").Value = -615
$ws.Range("N58").Value = -11775
# row 61
$ws.Range("H61").Value = 396.75
$ws.Range("I61").Value = 396.75
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1190.25
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1018.25
$ws.Range("N61").ClearContents()
# row 93
$ws.Range("H93").Value = 48000
$ws.Range("J93").Value = 48000
$ws.Range("L93").Value = 48000
$ws.Range("N93").Value = -52992
# row 96
$ws.Range("H96").Value = 2037.5
$ws.Range("I96").Value = 2271.4285
$ws.Range("J96").Value = 400
$ws.Range("K96").Value = 6814.2855
$ws.Range("L96").Value = 1200
$ws.Range("M96").Value = -5441.2855
$ws.Range("N96").Value = -3946
# row 116
$ws.Range("H116").Value = 6499
$ws.Range("I116").Value = 6499
$ws.Range("K116").Value = 6499
$ws.Range("M116").Value = -3057
# row 132
$ws.Range("H132").Value = 5245.5835
$ws.Range("I132").Value = 5533.222
$ws.Range("J132").Value = 4382.6665
$ws.Range("K132").Value = 16599.666
$ws.Range("L132").Value = 13147.9995
$ws.Range("M132").Value = -14069.666
$ws.Range("N132").Value = -18207.9995

$ws = $wb.Worksheets.Item("ARM")
# row 4
$ws.Range("H4").Value = 531.8
$ws.Range("I4").Value = 614.75
$ws.Range("K4").Value = 614.75
$ws.Range("M4").Value = -498.75
# row 5
$ws.Range("H5").Value = 261.1111
$ws.Range("I5").Value = 269
$ws.Range("K5").Value = 269
$ws.Range("M5").Value = -157
# row 32
$ws.Range("H32").Value = 7616.923
$ws.Range("I32").Value = 6793.375
$ws.Range("K32").Value = 6793.375
$ws.Range("M32").Value = -6506.375
# row 61
$ws.Range("H61").Value = 3126.3438
$ws.Range("I61").Value = 2418.7917
$ws.Range("K61").Value = 2418.7917
$ws.Range("M61").Value = -2206.7917
# row 95
$ws.Range("H95").Value = 31949
$ws.Range("J95").Value = 31949
$ws.Range("L95").Value = 31949
$ws.Range("N95").Value = -37441
# row 122
$ws.Range("H122").Value = 2461.8
$ws.Range("I122").Value = 2261
$ws.Range("K122").Value = 6783
$ws.Range("M122").Value = -4333
# row 132
$ws.Range("H132").Value = 2646.9
$ws.Range("J132").Value = 2599
$ws.Range("L132").Value = 7797
$ws.Range("N132").Value = -12857
# row 136
$ws.Range("H136").Value = 3126.3438
$ws.Range("I136").Value = 2418.7917
$ws.Range("K136").Value = 7256.375100000001
$ws.Range("M136").Value = -4706.375100000001

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 785.25
$ws.Range("I3").Value = 785.25
$ws.Range("K3").Value = 785.25
$ws.Range("M3").Value = -671.25
# row 4
$ws.Range("H4").Value = 261.1111
$ws.Range("I4").Value = 269
$ws.Range("K4").Value = 269
$ws.Range("M4").Value = -154
# row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
# row 20
$ws.Range("H20").Value = 1878.6364
$ws.Range("I20").Value = 2043.2941
$ws.Range("J20").Value = 1318.8
$ws.Range("K20").Value = 2043.2941
$ws.Range("L20").Value = 1318.8
$ws.Range("M20").Value = -1796.2941
$ws.Range("N20").Value = -1812.8
# row 99
$ws.Range("H99").Value = 3145
$ws.Range("J99").Value = 2996
$ws.Range("L99").Value = 2996
$ws.Range("N99").Value = -5992
# row 134
$ws.Range("H134").Value = 1297.6538
$ws.Range("I134").Value = 1285.56
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 3856.68
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -1321.68
$ws.Range("N134").Value = -9870

$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 2354650.2
$ws.Range("I22").Value = 1578.4445
$ws.Range("K22").Value = 1578.4445
$ws.Range("M22").Value = -1228.4445
# row 31
$ws.Range("H31").Value = 3876.75
$ws.Range("I31").Value = 2753.5
$ws.Range("K31").Value = 2753.5
$ws.Range("M31").Value = -2458.5
# row 34
$ws.Range("H34").Value = 3876.75
$ws.Range("I34").Value = 2753.5
$ws.Range("K34").Value = 2753.5
$ws.Range("M34").Value = -2551.5
# row 58
$ws.Range("H58").Value = 1731.875
$ws.Range("I58").Value = 1222.875
$ws.Range("K58").Value = 1222.875
$ws.Range("M58").Value = -1019.875
# row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# row 122
$ws.Range("H122").Value = 3059
$ws.Range("J122").Value = 3200
$ws.Range("L122").Value = 9600
$ws.Range("N122").Value = -14500
# row 134
$ws.Range("H134").Value = 1910.4166
$ws.Range("I134").Value = 1598.6666
$ws.Range("J134").Value = 2845.6667
$ws.Range("K134").Value = 4795.9998
$ws.Range("L134").Value = 8537.000100000001
$ws.Range("M134").Value = -2260.9998
$ws.Range("N134").Value = -13607.0001
# row 136
$ws.Range("H136").Value = 1731.875
$ws.Range("I136").Value = 1222.875
$ws.Range("K136").Value = 3668.625
$ws.Range("M136").Value = -1118.625

$ws = $wb.Worksheets.Item("CUL")
# row 7
$ws.Range("H7").Value = 80066.53999999999
$ws.Range("I7").Value = 250019.75
$ws.Range("J7").Value = 4531.778
$ws.Range("K7").Value = 750059.25
$ws.Range("L7").Value = 13595.334
$ws.Range("M7").Value = -749947.25
$ws.Range("N7").Value = -13819.334
# row 95
$ws.Range("H95").Value = 1124.5
$ws.Range("I95").Value = 1187
$ws.Range("J95").Value = 999.5
$ws.Range("K95").Value = 3561
$ws.Range("L95").Value = 2998.5
$ws.Range("M95").Value = -1502
$ws.Range("N95").Value = -7116.5
# row 97
$ws.Range("H97").Value = 1799
$ws.Range("I97").Value = 1799
$ws.Range("K97").Value = 5397
$ws.Range("M97").Value = -4901
# row 122
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1997
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 17973
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -22873
# row 131
$ws.Range("H131").Value = 1902.2
$ws.Range("I131").Value = 1750
$ws.Range("J131").Value = 1957.5454
$ws.Range("K131").Value = 5250
$ws.Range("L131").Value = 5872.6362
$ws.Range("M131").Value = -210
$ws.Range("N131").Value = -15952.6362

$ws = $wb.Worksheets.Item("GSM")
# row 46
$ws.Range("H46").Value = 25974.25
$ws.Range("J46").Value = 31299.334
$ws.Range("L46").Value = 31299.334
$ws.Range("N46").Value = -31611.334
# row 102
$ws.Range("H102").Value = 1950.25
$ws.Range("I102").Value = 1950.25
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1950.25
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -328.25
$ws.Range("N102").ClearContents()
# row 122
$ws.Range("H122").Value = 3665.6667
$ws.Range("I122").Value = 2498.5
$ws.Range("K122").Value = 7495.5
$ws.Range("M122").Value = -5045.5
# row 132
$ws.Range("H132").Value = 1820.5
$ws.Range("I132").Value = 1716.6538
$ws.Range("J132").Value = 2270.5
$ws.Range("K132").Value = 5149.9614
$ws.Range("L132").Value = 6811.5
$ws.Range("M132").Value = -2619.9614
$ws.Range("N132").Value = -11871.5

$ws = $wb.Worksheets.Item("LTW")
# row 93
$ws.Range("H93").Value = 3465.3333
$ws.Range("I93").Value = 3465.3333
$ws.Range("K93").Value = 3465.3333
$ws.Range("M93").Value = -2217.3333
# row 99
$ws.Range("H99").Value = 88500
$ws.Range("J99").Value = 88500
$ws.Range("L99").Value = 88500
$ws.Range("N99").Value = -94490
# row 101
$ws.Range("H101").Value = 20361
$ws.Range("J101").Value = 20361
$ws.Range("L101").Value = 20361
$ws.Range("N101").Value = -26851
# row 122
$ws.Range("H122").Value = 3460.7144
$ws.Range("I122").Value = 3399.5
$ws.Range("J122").Value = 3542.3333
$ws.Range("K122").Value = 10198.5
$ws.Range("L122").Value = 10626.9999
$ws.Range("M122").Value = -7748.5
$ws.Range("N122").Value = -15526.9999

$ws = $wb.Worksheets.Item("WVR")
# row 54
$ws.Range("H54").Value = 18401.334
$ws.Range("I54").Value = 1206
$ws.Range("J54").Value = 26999
$ws.Range("K54").Value = 1206
$ws.Range("L54").Value = 26999
$ws.Range("M54").Value = -686
$ws.Range("N54").Value = -28039
# row 136
$ws.Range("H136").Value = 13425.333
$ws.Range("I136").Value = 8910.532999999999
$ws.Range("J136").Value = 35999.332
$ws.Range("K136").Value = 26731.599
$ws.Range("L136").Value = 107997.996
$ws.Range("M136").Value = -24181.599
$ws.Range("N136").Value = -113097.996

